$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that ends the Java for-loop example:
#   "System.out.println(i);// 1,2,3,4,5,6,7,8,9,10"
# (it currently also hosts the _GoBack bookmark and is directly
#  followed by the paragraph that just contains "}")
# ------------------------------------------------------------------
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*1,2,3,4,5,6,7,8,9,10*") {
        $anchor = $i
        break
    }
}

$p = $d.Paragraphs($anchor)

# Create a placeholder paragraph right after it; InsertXML, applied to
# a collapsed range sitting on an (empty) placeholder paragraph, fills
# that paragraph in rather than appending a duplicate, so we can grow
# the placeholder into the whole block of new paragraphs below in one
# shot while keeping exact control of the OOXML (proofErr markers,
# the page-break hint, the bottom border, etc.).
$p.Range.InsertParagraphAfter()
$ph = $d.Paragraphs($anchor + 1)
$r = $ph.Range
$r.Collapse(0)

$w = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

$xml = "<w:p $w><w:r><w:t>}</w:t></w:r></w:p>" + `
       "<w:p $w/>" + `
       "<w:p $w/>" + `
       "<w:p $w><w:pPr><w:pBdr><w:bottom w:val=`"double`" w:sz=`"6`" w:space=`"1`" w:color=`"auto`"/></w:pBdr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>21-12-2021</w:t></w:r></w:p>" + `
       "<w:p $w><w:r><w:t xml:space=`"preserve`">Examples on </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>for  loop</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p>" + `
       "<w:p $w><w:r><w:t>Examples on nested for loops</w:t></w:r></w:p>" + `
       "<w:p $w><w:r><w:t>Find the factorial for given number</w:t></w:r></w:p>" + `
       "<w:p $w><w:r><w:t>WAP to achieve below output</w:t></w:r></w:p>" + `
       "<w:p $w/>"

$r.InsertXML($xml)

# ------------------------------------------------------------------
# Move the (single, special) _GoBack bookmark from its old home -
# the end of the "1,2,3,4,5,6,7,8,9,10" paragraph - onto the new
# empty paragraph that now follows "WAP to achieve below output".
# Bookmarks.Add with the existing bookmark's name relocates it.
# ------------------------------------------------------------------
$bookmarkHost = $d.Paragraphs($anchor + 9)
$d.Bookmarks.Add("_GoBack", $bookmarkHost.Range)

# ------------------------------------------------------------------
# The paragraph that used to hold the trailing "}" (now pushed one
# further down, right after the bookmark paragraph) loses its text
# and gets five spaces instead.
# ------------------------------------------------------------------
$tail = $d.Paragraphs($anchor + 10)
$tail.Range.Text = "     "
